# Auto-generated update of leve profit columns (H:N) across all 8 item-category sheets.
# Mirrors a scheduled price-refresh run: currentAveragePrice* and LeveProfit* columns
# are recomputed from freshly pulled market data; cells with no resulting value are cleared
# rather than written as 0/blank placeholders (matching the sparse source formatting).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Cells.Item(2, 8).Value = 103
$ws.Cells.Item(2, 9).Value = 81.25
$ws.Cells.Item(2, 10).Value = 124.75
$ws.Cells.Item(2, 11).Value = 81.25
$ws.Cells.Item(2, 12).Value = 124.75
$ws.Cells.Item(2, 13).Value = 31.75
$ws.Cells.Item(2, 14).Value = -350.75
# Row 5 (Leve Item ID 5503)
$ws.Cells.Item(5, 8).Value = 150
$ws.Cells.Item(5, 9).Value = 135
$ws.Cells.Item(5, 10).Value = 172.5
$ws.Cells.Item(5, 11).Value = 135
$ws.Cells.Item(5, 12).Value = 172.5
$ws.Cells.Item(5, 13).Value = -20
$ws.Cells.Item(5, 14).Value = -402.5
# Row 7 (Leve Item ID 1960)
$ws.Cells.Item(7, 8).Value = 4371.25
$ws.Cells.Item(7, 9).Value = 3952.5
$ws.Cells.Item(7, 10).Value = 4790
$ws.Cells.Item(7, 11).Value = 3952.5
$ws.Cells.Item(7, 12).Value = 4790
$ws.Cells.Item(7, 13).Value = -3840.5
$ws.Cells.Item(7, 14).Value = -5014
# Row 8 (Leve Item ID 4565)
$ws.Cells.Item(8, 8).Value = 120.72727
$ws.Cells.Item(8, 9).Value = 122.8
$ws.Cells.Item(8, 11).Value = 368.4
$ws.Cells.Item(8, 13).Value = -229.4
# Row 14 (Leve Item ID 1960)
$ws.Cells.Item(14, 8).Value = 4371.25
$ws.Cells.Item(14, 9).Value = 3952.5
$ws.Cells.Item(14, 10).Value = 4790
$ws.Cells.Item(14, 11).Value = 3952.5
$ws.Cells.Item(14, 12).Value = 4790
$ws.Cells.Item(14, 13).Value = -3761.5
$ws.Cells.Item(14, 14).Value = -5172
# Row 17 (Leve Item ID 38956)
$ws.Cells.Item(17, 8).Value = 2780.6667
$ws.Cells.Item(17, 10).Value = 2780.6667
$ws.Cells.Item(17, 12).Value = 8342.000100000001
$ws.Cells.Item(17, 14).Value = -8678.000100000001
# Row 41 (Leve Item ID 5478)
$ws.Cells.Item(41, 8).Value = 1467.3334
$ws.Cells.Item(41, 9).Value = 2001
$ws.Cells.Item(41, 10).Value = 400
$ws.Cells.Item(41, 11).Value = 2001
$ws.Cells.Item(41, 12).Value = 400
$ws.Cells.Item(41, 13).Value = -1561
$ws.Cells.Item(41, 14).Value = -1280
# Row 87 (Leve Item ID 10651)
$ws.Cells.Item(87, 8).Value = 49999
$ws.Cells.Item(87, 9).Value = 49999
$ws.Cells.Item(87, 11).Value = 49999
$ws.Cells.Item(87, 13).Value = -48751
# Row 90 (Leve Item ID 10651)
$ws.Cells.Item(90, 8).Value = 49999
$ws.Cells.Item(90, 9).Value = 49999
$ws.Cells.Item(90, 11).Value = 149997
$ws.Cells.Item(90, 13).Value = -143757
# Row 116 (Leve Item ID 27778)
$ws.Cells.Item(116, 8).Value = 14880
$ws.Cells.Item(116, 9).Value = 3900
$ws.Cells.Item(116, 10).Value = 17625
$ws.Cells.Item(116, 11).Value = 3900
$ws.Cells.Item(116, 12).Value = 17625
$ws.Cells.Item(116, 13).Value = -458
$ws.Cells.Item(116, 14).Value = -24509
# Row 132 (Leve Item ID 44049)
$ws.Cells.Item(132, 8).Value = 9611.799999999999
$ws.Cells.Item(132, 9).Value = 9611.799999999999
$ws.Cells.Item(132, 11).Value = 28835.4
$ws.Cells.Item(132, 13).Value = -26305.4

$ws = $wb.Worksheets.Item("ARM")
# Row 3 (Leve Item ID 2494)
$ws.Cells.Item(3, 8).Value = 187.5
$ws.Cells.Item(3, 9).Value = 187.5
$ws.Cells.Item(3, 11).Value = 187.5
$ws.Cells.Item(3, 13).Value = -72.5
# Row 13 (Leve Item ID 2656)
$ws.Cells.Item(13, 8).Value = 999
$ws.Cells.Item(13, 10).Value = 999
$ws.Cells.Item(13, 12).Value = 999
$ws.Cells.Item(13, 14).Value = -1287
# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 5971.3335
$ws.Cells.Item(32, 9).Value = 3592.75
$ws.Cells.Item(32, 11).Value = 3592.75
$ws.Cells.Item(32, 13).Value = -3305.75
# Row 45 (Leve Item ID 27714)
$ws.Cells.Item(45, 8).Value = 1496.8
$ws.Cells.Item(45, 9).Value = 1498.5
$ws.Cells.Item(45, 10).Value = 1490
$ws.Cells.Item(45, 11).Value = 1498.5
$ws.Cells.Item(45, 12).Value = 1490
$ws.Cells.Item(45, 13).Value = -1121.5
$ws.Cells.Item(45, 14).Value = -2244
# Row 62 (Leve Item ID 10719)
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 14).ClearContents()
# Row 65 (Leve Item ID 10719)
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 14).ClearContents()
# Row 74 (Leve Item ID 44000)
$ws.Cells.Item(74, 8).Value = 5247.5
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()
# Row 77 (Leve Item ID 44000)
$ws.Cells.Item(77, 8).Value = 5247.5
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()
# Row 97 (Leve Item ID 19941)
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).ClearContents()
$ws.Cells.Item(97, 14).ClearContents()
# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 6500
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 6500
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 19500
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).Value = -24400

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID 19939)
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).ClearContents()
# Row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 8003.6665
$ws.Cells.Item(134, 9).Value = 2012
$ws.Cells.Item(134, 10).Value = 10999.5
$ws.Cells.Item(134, 11).Value = 6036
$ws.Cells.Item(134, 12).Value = 32998.5
$ws.Cells.Item(134, 13).Value = -3501
$ws.Cells.Item(134, 14).Value = -38068.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Cells.Item(16, 8).Value = 2197.5
$ws.Cells.Item(16, 9).Value = 2316
$ws.Cells.Item(16, 11).Value = 2316
$ws.Cells.Item(16, 13).Value = -2029
# Row 58 (Leve Item ID 44021)
$ws.Cells.Item(58, 8).Value = 25000
$ws.Cells.Item(58, 10).Value = 25000
$ws.Cells.Item(58, 12).Value = 25000
$ws.Cells.Item(58, 14).Value = -25406
# Row 113 (Leve Item ID 27691)
$ws.Cells.Item(113, 8).Value = 2197.5
$ws.Cells.Item(113, 9).Value = 2316
$ws.Cells.Item(113, 11).Value = 2316
$ws.Cells.Item(113, 13).Value = -146
# Row 134 (Leve Item ID 44020)
$ws.Cells.Item(134, 8).Value = 4197
$ws.Cells.Item(134, 9).Value = 4658
$ws.Cells.Item(134, 11).Value = 13974
$ws.Cells.Item(134, 13).Value = -11439
# Row 136 (Leve Item ID 44021)
$ws.Cells.Item(136, 8).Value = 25000
$ws.Cells.Item(136, 10).Value = 25000
$ws.Cells.Item(136, 12).Value = 75000
$ws.Cells.Item(136, 14).Value = -80100

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Leve Item ID 4847)
$ws.Cells.Item(2, 8).Value = 98.5
$ws.Cells.Item(2, 9).Value = 101.25
$ws.Cells.Item(2, 10).Value = 95.75
$ws.Cells.Item(2, 11).Value = 607.5
$ws.Cells.Item(2, 12).Value = 574.5
$ws.Cells.Item(2, 13).Value = -494.5
$ws.Cells.Item(2, 14).Value = -800.5
# Row 86 (Leve Item ID 12892)
$ws.Cells.Item(86, 8).Value = 503
$ws.Cells.Item(86, 10).Value = 503
$ws.Cells.Item(86, 12).Value = 1509
$ws.Cells.Item(86, 14).Value = -3881
# Row 89 (Leve Item ID 12892)
$ws.Cells.Item(89, 8).Value = 503
$ws.Cells.Item(89, 10).Value = 503
$ws.Cells.Item(89, 12).Value = 4527
$ws.Cells.Item(89, 14).Value = -16383
# Row 102 (Leve Item ID 19813)
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Cells.Item(122, 8).Value = 2499
$ws.Cells.Item(122, 9).Value = 2499
$ws.Cells.Item(122, 11).Value = 7497
$ws.Cells.Item(122, 13).Value = -5047
# Row 126 (Leve Item ID 36184)
$ws.Cells.Item(126, 8).Value = 6338.6665
$ws.Cells.Item(126, 10).Value = 2994
$ws.Cells.Item(126, 12).Value = 8982
$ws.Cells.Item(126, 14).Value = -13922
# Row 128 (Leve Item ID 34544)
$ws.Cells.Item(128, 8).Value = 60000
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 60000
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 60000
$ws.Cells.Item(128, 13).ClearContents()
$ws.Cells.Item(128, 14).Value = -69960

$ws = $wb.Worksheets.Item("LTW")
# Row 96 (Leve Item ID 19735)
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 14).ClearContents()
# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 5987.8
$ws.Cells.Item(132, 9).Value = 5987.8
$ws.Cells.Item(132, 11).Value = 17963.4
$ws.Cells.Item(132, 13).Value = -15433.4

$ws = $wb.Worksheets.Item("WVR")
# Row 68 (Leve Item ID 10762)
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
# Row 69 (Leve Item ID 10951)
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()
# Row 71 (Leve Item ID 10762)
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
# Row 72 (Leve Item ID 10951)
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()
# Row 96 (Leve Item ID 19977)
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 13).ClearContents()
# Row 107 (Leve Item ID 27746)
$ws.Cells.Item(107, 8).Value = 746.8570999999999
$ws.Cells.Item(107, 9).Value = 746.8570999999999
$ws.Cells.Item(107, 11).Value = 2240.5713
$ws.Cells.Item(107, 13).Value = -320.5712999999996
# Row 126 (Leve Item ID 36210)
$ws.Cells.Item(126, 8).Value = 2440.8
$ws.Cells.Item(126, 9).Value = 2101.3333
$ws.Cells.Item(126, 10).Value = 2950
$ws.Cells.Item(126, 11).Value = 6303.999899999999
$ws.Cells.Item(126, 12).Value = 8850
$ws.Cells.Item(126, 13).Value = -3833.999899999999
$ws.Cells.Item(126, 14).Value = -13790
# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, 8).Value = 1043.4445
$ws.Cells.Item(132, 9).Value = 1006.7143
$ws.Cells.Item(132, 11).Value = 3020.1429
$ws.Cells.Item(132, 13).Value = -490.1428999999998
# Row 136 (Leve Item ID 44031)
$ws.Cells.Item(136, 8).Value = 3472.7827
$ws.Cells.Item(136, 9).Value = 3176.0908
$ws.Cells.Item(136, 11).Value = 9528.2724
$ws.Cells.Item(136, 13).Value = -6978.2724
